$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 511.44446
$ws.Range("I29").Value = 340.6
$ws.Range("J29").Value = 725
$ws.Range("K29").Value = 1021.8
$ws.Range("L29").Value = 2175
$ws.Range("M29").Value = -740.8000000000001
$ws.Range("N29").Value = -2737
$ws.Range("H38").Value = 1898182.2
$ws.Range("J38").Value = 2425
$ws.Range("L38").Value = 7275
$ws.Range("N38").Value = -8019
$ws.Range("H43").Value = 2564.4443
$ws.Range("I43").Value = 2642.2593
$ws.Range("J43").Value = 2331
$ws.Range("K43").Value = 2642.2593
$ws.Range("L43").Value = 2331
$ws.Range("M43").Value = -2573.2593
$ws.Range("N43").Value = -2469
$ws.Range("H58").Value = 851634.6
$ws.Range("I58").Value = 1623620.6
$ws.Range("J58").Value = 2450
$ws.Range("K58").Value = 4870861.800000001
$ws.Range("L58").Value = 7350
$ws.Range("M58").Value = -4870711.800000001
$ws.Range("N58").Value = -7650
$ws.Range("H61").Value = 2211187.5
$ws.Range("I61").Value = 4762040.5
$ws.Range("J61").Value = 24741.857
$ws.Range("K61").Value = 14286121.5
$ws.Range("L61").Value = 74225.571
$ws.Range("M61").Value = -14285949.5
$ws.Range("N61").Value = -74569.571
$ws.Range("H87").Value = 31838
$ws.Range("J87").Value = 31838
$ws.Range("L87").Value = 31838
$ws.Range("N87").Value = -34334
$ws.Range("H90").Value = 31838
$ws.Range("J90").Value = 31838
$ws.Range("L90").Value = 95514
$ws.Range("N90").Value = -107994
$ws.Range("H112").Value = 942.4783
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 958.0454999999999
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 2874.1365
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -5090.1365
$ws.Range("H129").Value = 957.0769
$ws.Range("I129").Value = 575.1111
$ws.Range("J129").Value = 1071.6666
$ws.Range("K129").Value = 1725.3333
$ws.Range("L129").Value = 3214.9998
$ws.Range("M129").Value = 3274.6667
$ws.Range("N129").Value = -13214.9998
$ws.Range("H132").Value = 4468809.5
$ws.Range("I132").Value = 4812160
$ws.Range("J132").Value = 5251.5
$ws.Range("K132").Value = 14436480
$ws.Range("L132").Value = 15754.5
$ws.Range("M132").Value = -14433950
$ws.Range("N132").Value = -20814.5
$ws.Range("H135").Value = 1428.32
$ws.Range("I135").Value = 865.1
$ws.Range("J135").Value = 3681.2
$ws.Range("K135").Value = 7785.900000000001
$ws.Range("L135").Value = 33130.8
$ws.Range("M135").Value = -5250.900000000001
$ws.Range("N135").Value = -38200.8
$ws.Range("H138").Value = 6821.3335
$ws.Range("I138").Value = 1332.9375
$ws.Range("J138").Value = 17798.125
$ws.Range("K138").Value = 3998.8125
$ws.Range("L138").Value = 53394.375
$ws.Range("M138").Value = 1141.1875
$ws.Range("N138").Value = -63674.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2139.5144
$ws.Range("I61").Value = 1349.0555
$ws.Range("J61").Value = 2976.4707
$ws.Range("K61").Value = 1349.0555
$ws.Range("L61").Value = 2976.4707
$ws.Range("M61").Value = -1137.0555
$ws.Range("N61").Value = -3400.4707
$ws.Range("H132").Value = 2740.2
$ws.Range("I132").Value = 2589.4595
$ws.Range("J132").Value = 4599.3335
$ws.Range("K132").Value = 7768.3785
$ws.Range("L132").Value = 13798.0005
$ws.Range("M132").Value = -5238.3785
$ws.Range("N132").Value = -18858.0005
$ws.Range("H136").Value = 2139.5144
$ws.Range("I136").Value = 1349.0555
$ws.Range("J136").Value = 2976.4707
$ws.Range("K136").Value = 4047.1665
$ws.Range("L136").Value = 8929.4121
$ws.Range("M136").Value = -1497.1665
$ws.Range("N136").Value = -14029.4121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2687.8298
$ws.Range("I134").Value = 2571.1052
$ws.Range("K134").Value = 7713.3156
$ws.Range("M134").Value = -5178.3156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 42523.54
$ws.Range("I31").Value = 1361.2
$ws.Range("J31").Value = 56717.45
$ws.Range("K31").Value = 1361.2
$ws.Range("L31").Value = 56717.45
$ws.Range("M31").Value = -1066.2
$ws.Range("N31").Value = -57307.45
$ws.Range("H34").Value = 42523.54
$ws.Range("I34").Value = 1361.2
$ws.Range("J34").Value = 56717.45
$ws.Range("K34").Value = 1361.2
$ws.Range("L34").Value = 56717.45
$ws.Range("M34").Value = -1159.2
$ws.Range("N34").Value = -57121.45
$ws.Range("H107").Value = 5078.9546
$ws.Range("I107").Value = 8290.615
$ws.Range("J107").Value = 439.8889
$ws.Range("K107").Value = 8290.615
$ws.Range("L107").Value = 439.8889
$ws.Range("M107").Value = -6370.615
$ws.Range("N107").Value = -4279.8889
$ws.Range("H132").Value = 22729862
$ws.Range("I132").Value = 20002284
$ws.Range("K132").Value = 60006852
$ws.Range("M132").Value = -60004322
$ws.Range("H134").Value = 1045
$ws.Range("I134").Value = 925.26086
$ws.Range("J134").Value = 1389.25
$ws.Range("K134").Value = 2775.78258
$ws.Range("L134").Value = 4167.75
$ws.Range("M134").Value = -240.7825800000001
$ws.Range("N134").Value = -9237.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 109921.63
$ws.Range("I70").Value = 226089.11
$ws.Range("J70").Value = 5370.9
$ws.Range("K70").Value = 226089.11
$ws.Range("L70").Value = 5370.9
$ws.Range("M70").Value = -225819.11
$ws.Range("N70").Value = -5910.9
$ws.Range("H73").Value = 109921.63
$ws.Range("I73").Value = 226089.11
$ws.Range("J73").Value = 5370.9
$ws.Range("K73").Value = 226089.11
$ws.Range("L73").Value = 5370.9
$ws.Range("M73").Value = -225153.11
$ws.Range("N73").Value = -7242.9
$ws.Range("H80").Value = 125130296
$ws.Range("I80").Value = 166839730
$ws.Range("K80").Value = 166839730
$ws.Range("M80").Value = -166838732
$ws.Range("H83").Value = 125130296
$ws.Range("I83").Value = 166839730
$ws.Range("K83").Value = 834198650
$ws.Range("M83").Value = -834193658
$ws.Range("H132").Value = 2435.568
$ws.Range("I132").Value = 1376.2258
$ws.Range("J132").Value = 4961.6924
$ws.Range("K132").Value = 4128.6774
$ws.Range("L132").Value = 14885.0772
$ws.Range("M132").Value = -1598.6774
$ws.Range("N132").Value = -19945.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1958.8
$ws.Range("I136").Value = 1750.9231
$ws.Range("J136").Value = 3310
$ws.Range("K136").Value = 5252.7693
$ws.Range("L136").Value = 9930
$ws.Range("M136").Value = -2702.7693
$ws.Range("N136").Value = -15030

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 57266.668
$ws.Range("J123").Value = 57266.668
$ws.Range("L123").Value = 57266.668
$ws.Range("N123").Value = -67066.66800000001
$ws.Range("H132").Value = 1632.1184
$ws.Range("I132").Value = 1496.6615
$ws.Range("J132").Value = 2432.5454
$ws.Range("K132").Value = 4489.9845
$ws.Range("L132").Value = 7297.6362
$ws.Range("M132").Value = -1959.9845
$ws.Range("N132").Value = -12357.6362
